$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cells E16:E18 - change value from "NA" to "removeNegative" and clear the style
$ws.Range("E16").Value = "removeNegative"
$ws.Range("E17").Value = "removeNegative"
$ws.Range("E18").Value = "removeNegative"

$ws.Range("E16:E18").Style = "Normal"

# Update the active selection to E18
$ws.Range("E18").Select()
